# Append a new "CHP" row (element=CHP1 / type=CHP) right after the
# existing "bat1"/"bat" row, extending the table from A1:B5 to A1:B6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "CHP1"
$ws.Range("B6").Value = "CHP"
